$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new transaction data
$ws.Range("B2").Value = "2023-04-17T12:34"
$ws.Range("C2").Value = "Master"
$ws.Range("D2").Value = "Cash"
$ws.Range("E2").Value = "advance"
$ws.Range("F2").Value = 44999

# Remove rows 3-5 (the other transactions no longer needed)
$ws.Range("A3:F5").Delete()
